$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: capture old row 25 data (IGFBP3 "Data from..." row) before shifting ---
$origRow25 = $ws.Range("A25:Q25").Value()

# --- Step 2: insert a new row at position 12 (shifts old rows 12-25 down to 13-26) ---
$ws.Rows.Item(12).Insert()

# --- Step 3: populate new row 12 with the captured row-25 data, forcing text (no date/number coercion) ---
$rng12 = $ws.Range("A12:Q12")
$rng12.NumberFormat = "@"
$rng12.Value = $origRow25
$rng12.ClearFormats()

# --- Step 4: remove the now-duplicated old row (shifted to row 26) ---
$ws.Rows.Item(26).Delete()

# --- Step 5: append four brand-new bibliography rows (26-29) ---
$arr26 = New-Object 'object[,]' 1,17
$arr26[0,0] = 'Erin Muhlbradt, Ekaterina Asatiani, Elizabeth Ortner, Antai Wang, Edward P. Gelmann'
$arr26[0,1] = '; ; ; ; '
$arr26[0,2] = 'https://openalex.org/W4392679220'
$arr26[0,3] = 'Data from NKX3.1 Activates Expression of Insulin-like Growth Factor Binding Protein-3 to Mediate Insulin-like Growth Factor-I Signaling and Cell Proliferation'
$arr26[0,4] = '2023-03-30'
$arr26[0,5] = 'N/A'
$arr26[0,6] = 'N/A'
$arr26[0,7] = 'https://doi.org/10.1158/0008-5472.c.6499377'
$arr26[0,8] = 'N/A'
$arr26[0,9] = 'submittedVersion'
$arr26[0,10] = 'closed'
$arr26[0,11] = 'en'
$arr26[0,12] = '0'
$arr26[0,13] = '2023'
$arr26[0,14] = 'NA'
$arr26[0,15] = 'https://doi.org/10.1158/0008-5472.c.6499377'
$arr26[0,16] = 'article'
$rng26 = $ws.Range("A26:Q26")
$rng26.NumberFormat = "@"
$rng26.Value = $arr26
$rng26.ClearFormats()

$arr27 = New-Object 'object[,]' 1,17
$arr27[0,0] = 'Bin Cai, Zheng Tian, Edward P. Gelmann'
$arr27[0,1] = '; ; '
$arr27[0,2] = 'https://openalex.org/W4392679801'
$arr27[0,3] = 'Data from NKX3.1 Suppresses &lt;i&gt;TMPRSS2–ERG&lt;/i&gt; Gene Rearrangement and Mediates Repair of Androgen Receptor–Induced DNA Damage'
$arr27[0,4] = '2023-03-30'
$arr27[0,5] = 'N/A'
$arr27[0,6] = 'N/A'
$arr27[0,7] = 'https://doi.org/10.1158/0008-5472.c.6506729'
$arr27[0,8] = 'N/A'
$arr27[0,9] = 'submittedVersion'
$arr27[0,10] = 'closed'
$arr27[0,11] = 'en'
$arr27[0,12] = '0'
$arr27[0,13] = '2023'
$arr27[0,14] = 'NA'
$arr27[0,15] = 'https://doi.org/10.1158/0008-5472.c.6506729'
$arr27[0,16] = 'article'
$rng27 = $ws.Range("A27:Q27")
$rng27.NumberFormat = "@"
$rng27.Value = $arr27
$rng27.ClearFormats()

$arr28 = New-Object 'object[,]' 1,17
$arr28[0,0] = 'Bin Cai, Zheng Tian, Edward P. Gelmann'
$arr28[0,1] = '; ; '
$arr28[0,2] = 'https://openalex.org/W4392687152'
$arr28[0,3] = 'Data from NKX3.1 Suppresses &lt;i&gt;TMPRSS2–ERG&lt;/i&gt; Gene Rearrangement and Mediates Repair of Androgen Receptor–Induced DNA Damage'
$arr28[0,4] = '2023-03-30'
$arr28[0,5] = 'N/A'
$arr28[0,6] = 'N/A'
$arr28[0,7] = 'https://doi.org/10.1158/0008-5472.c.6506729.v1'
$arr28[0,8] = 'N/A'
$arr28[0,9] = 'submittedVersion'
$arr28[0,10] = 'closed'
$arr28[0,11] = 'en'
$arr28[0,12] = '0'
$arr28[0,13] = '2023'
$arr28[0,14] = 'NA'
$arr28[0,15] = 'https://doi.org/10.1158/0008-5472.c.6506729.v1'
$arr28[0,16] = 'article'
$rng28 = $ws.Range("A28:Q28")
$rng28.NumberFormat = "@"
$rng28.Value = $arr28
$rng28.ClearFormats()

$arr29 = New-Object 'object[,]' 1,17
$arr29[0,0] = 'Shyh‐Han Tan, Ayush Dagvadorj, Feng Shen, Lin Gu, Zhiyong Liao, Junaid Abdulghani, Ying Zhang, Edward P. Gelmann, Tobias Zellweger, Zoran Čulig, Tapio Visakorpi, Lukas Bubendorf, Robert A. Kirken, James G. Karras, Marja T. Nevalainen'
$arr29[0,1] = '; ; ; ; ; ; ; ; ; ; ; ; ; ; '
$arr29[0,2] = 'https://openalex.org/W4392692487'
$arr29[0,3] = 'Data from Transcription Factor Stat5 Synergizes with Androgen Receptor in Prostate Cancer Cells'
$arr29[0,4] = '2023-03-30'
$arr29[0,5] = 'N/A'
$arr29[0,6] = 'N/A'
$arr29[0,7] = 'https://doi.org/10.1158/0008-5472.c.6496695.v1'
$arr29[0,8] = 'N/A'
$arr29[0,9] = 'submittedVersion'
$arr29[0,10] = 'closed'
$arr29[0,11] = 'en'
$arr29[0,12] = '0'
$arr29[0,13] = '2023'
$arr29[0,14] = 'NA'
$arr29[0,15] = 'https://doi.org/10.1158/0008-5472.c.6496695.v1'
$arr29[0,16] = 'article'
$rng29 = $ws.Range("A29:Q29")
$rng29.NumberFormat = "@"
$rng29.Value = $arr29
$rng29.ClearFormats()
